$d = $word.ActiveDocument

# 1. Remove the stray _GoBack bookmark near the title; bookmark ids 2-7 renumber to 1-6
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 2. Fix the ".resilient/app.config" reference to "app.config" with proofing marks
$appConfigPara = $d.Paragraphs(92)
$appConfigXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText"/><w:ind w:left="360"/></w:pPr><w:r><w:t>The l</w:t></w:r><w:r w:rsidR="00623A24"><w:t xml:space="preserve">og is controlled in the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00623A24" w:rsidRPr="002965D4"><w:rPr><w:rStyle w:val="codeChar"/></w:rPr><w:t>app.config</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00623A24"><w:t xml:space="preserve"> file under the section </w:t></w:r><w:r w:rsidR="00623A24" w:rsidRPr="008A31F5"><w:rPr><w:rStyle w:val="codeChar"/></w:rPr><w:t>[resil</w:t></w:r><w:r w:rsidR="002965D4" w:rsidRPr="008A31F5"><w:rPr><w:rStyle w:val="codeChar"/></w:rPr><w:t>i</w:t></w:r><w:r w:rsidR="00623A24" w:rsidRPr="008A31F5"><w:rPr><w:rStyle w:val="codeChar"/></w:rPr><w:t>ent]</w:t></w:r><w:r w:rsidR="00623A24"><w:t xml:space="preserve"> and the property </w:t></w:r><w:r w:rsidR="00623A24" w:rsidRPr="002965D4"><w:rPr><w:rStyle w:val="codeChar"/></w:rPr><w:t>logdir</w:t></w:r><w:r w:rsidR="00623A24"><w:t xml:space="preserve">. The default file name is </w:t></w:r><w:r w:rsidR="00623A24" w:rsidRPr="002965D4"><w:rPr><w:rStyle w:val="codeChar"/></w:rPr><w:t>app.log</w:t></w:r><w:r w:rsidR="00623A24"><w:t xml:space="preserve">. Each function will </w:t></w:r><w:r w:rsidR="008A31F5"><w:t>create</w:t></w:r><w:r w:rsidR="00623A24"><w:t xml:space="preserve"> progress information. Failures will show up </w:t></w:r><w:r w:rsidR="008A31F5"><w:t xml:space="preserve">as errors and may contain </w:t></w:r><w:r w:rsidR="00623A24"><w:t>python trace statements.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$appConfigPara.Range.InsertXML($appConfigXml)

# 3. Rewrite the Support paragraph text/hyperlink (and drop keepNext) and add a fresh _GoBack bookmark at its end
$supportPara = $d.Paragraphs(94)
$supportXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Please review the resilient-circuits log file to help identify your issue. </w:t></w:r><w:r><w:t xml:space="preserve">For additional support, </w:t></w:r><w:r><w:t xml:space="preserve">refer to the IBM Resilient Community forum: </w:t></w:r><w:hyperlink r:id="rId12" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>https://ibm.biz/resilientcommunity</w:t></w:r></w:hyperlink><w:r><w:t>.</w:t></w:r><w:bookmarkStart w:id="7" w:name="_GoBack"/><w:bookmarkEnd w:id="7"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$supportPara.Range.InsertXML($supportXml)

# 4. Delete the trailing "Including relevant information..." paragraph entirely
$extraPara = $d.Paragraphs(95)
$extraPara.Range.Delete()

Write-Host "edit complete"
